{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst oldText = \"Make sure ghost are not shown until \\u201cPlayer One\\u201d text disappears.\";\nconst newText = \"Add player and ghost collision detection.\";\n\nlet firstIdx = -1;\nlet dupIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (firstIdx === -1 && t === oldText) {\n    firstIdx = i;\n  } else if (dupIdx === -1 && t === newText) {\n    dupIdx = i;\n  }\n}\n\nif (firstIdx !== -1) {\n  // Replace the text of the first paragraph's run while keeping the\n  // paragraph (and its bookmark) in place.\n  const range = paragraphs.items[firstIdx].getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nif (dupIdx !== -1) {\n  // Remove the now-duplicate paragraph entirely.\n  paragraphs.items[dupIdx].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"Make sure ghost are not shown until \" + [char]0x201C + \"Player One\" + [char]0x201D + \" text disappears.\"\n$newText = \"Add player and ghost collision detection.\"\n\n# Replace the text of the first paragraph (keeps the paragraph, and the\n# bookmark it carries, in place) rather than touching the whole range so the\n# bookmark/paragraph mark survive intact.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n$find.Execute([ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$newText, [ref]1)\n\n# Now remove the now-duplicate paragraph (\"Add player and ghost collision\n# detection.\") that used to follow it.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq $newText) {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
